# Generate Report for Handback
# Updates the "68b60577-ecc8-423a-96c7-45905dd27208" row (row 7) on both the
# zh-cn and de-de sheets to reflect a new handback attempt that failed
# because the handback file version was not the latest one.

$wb = $excel.ActiveWorkbook

$fileMd       = "68b60577-ecc8-423a-96c7-45905dd27208.md"
$handoffXlfZh = "68b60577-ecc8-423a-96c7-45905dd27208.1eba7e1abf6c9b0faa775fc2ab7a06d819e71735.zh-cn.xlf"
$handoffXlfDe = "68b60577-ecc8-423a-96c7-45905dd27208.1eba7e1abf6c9b0faa775fc2ab7a06d819e71735.de-de.xlf"
$handbackTimeZh = "2016-08-29 17:01:35"
$handbackTimeDe = "2016-08-29 17:01:42"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7330fcd9b00dbd8d6f2a621c3fd03fa10e3de2e6/e2e/68b60577-ecc8-423a-96c7-45905dd27208.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d05389d5885c3f308feb28a6247f2a4b7e91790/e2e/68b60577-ecc8-423a-96c7-45905dd27208.md."
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d05389d5885c3f308feb28a6247f2a4b7e91790/e2e/68b60577-ecc8-423a-96c7-45905dd27208.md"

function Update-HandbackRow {
    param($ws, $handoffXlf, $handbackTime)

    # Latest Target File (I7) now has a value and becomes a hyperlink to the
    # source handback markdown file, just like the other rows in the table.
    $ws.Range("I7").Value = $fileMd
    $ws.Hyperlinks.Add($ws.Range("I7"), $hyperlinkTarget, $null, $null, $fileMd) | Out-Null
    $ws.Range("I7").Style = "HyperLink"
    $ws.Range("I7").Font.Underline = 2
    $ws.Range("I7").Font.Color = 15570276
    $ws.Range("I7").Font.Name = "Calibri"
    $ws.Range("I7").Font.Size = 11

    # Latest Handback File (J7)
    $ws.Range("J7").Value = $handoffXlf

    # Latest Handback DateTime (K7)
    $ws.Range("K7").Value = $handbackTime

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn $handoffXlfZh $handbackTimeZh

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe $handoffXlfDe $handbackTimeDe
